$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11-13 down to 12-14
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new weekly data point
$ws.Cells.Item(11, 1).Value = 8
$ws.Cells.Item(11, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(11, 3).Value = "Coquimbo"
$ws.Cells.Item(11, 4).Value = 44784
$ws.Cells.Item(11, 5).Value = 4
$ws.Cells.Item(11, 6).Value = 100112013
$ws.Cells.Item(11, 7).Value = "Alcachofa"
$ws.Cells.Item(11, 8).Value = "Madrigal"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 520
$ws.Cells.Item(11, 11).Value = 11500
$ws.Cells.Item(11, 12).Value = 12000
$ws.Cells.Item(11, 13).Value = 11750
$ws.Cells.Item(11, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(11, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(11, 16).Value = 294
$ws.Cells.Item(11, 17).Value = 40
$ws.Cells.Item(11, 18).Value = "Hortaliza"
